$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "Place:  Quex Park, Birchington, Kent" -> "Place:  Quex House, Birchington, Kent"
#    The single run " Park, Birchington, Kent" becomes two runs:
#    " House" and ", Birchington, Kent" (identical formatting).
# -----------------------------------------------------------------
$para1 = $d.Paragraphs(2).Range
$find1 = $para1.Find
$find1.Execute(" Park, Birchington, Kent", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1 = $find1.Parent
$start1 = $r1.Start
$sub1a = $d.Range($start1, $start1 + 5)
$sub1a.Text = " House"
$sub1b = $d.Range($start1 + 6, $start1 + 25)
# Toggle Bold off/on (net no-op on formatting) to force Word to split this
# into its own run without altering any character formatting.
$sub1b.Bold = 0
$sub1b.Bold = 1

# -----------------------------------------------------------------
# 2) "Event: Two intrepid women from Quex Park, Kent travelled..." ->
#    "...Quex House, Kent travelled..."
#    The single run " Park," becomes two runs: " House" and ",".
# -----------------------------------------------------------------
$para2 = $d.Paragraphs(3).Range
$find2 = $para2.Find
$find2.Execute(" Park,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2 = $find2.Parent
$start2 = $r2.Start
$sub2a = $d.Range($start2, $start2 + 5)
$sub2a.Text = " House"
$sub2b = $d.Range($start2 + 6, $start2 + 7)
$sub2b.Bold = 0
$sub2b.Bold = 1

# -----------------------------------------------------------------
# 3) Endnote 1: merge the three runs (split by gramStart/gramEnd
#    proofErr markers around "To") into a single run with the
#    unbroken text.
# -----------------------------------------------------------------
$endnoteText = "Hutchinson, Pamela (2019) Where To Begin With Early Women Filmmakers p. BFI"
$en1 = $d.Endnotes(1).Range
$en1.Text = $endnoteText + " "
$en1b = $d.Endnotes(1).Range
$en1b.Text = $endnoteText
